# fayoumi (add course module)
#
# Updates the "register" sample sheet:
#  - A2 name text revised
#  - B2 phone number re-entered as text (quote-prefixed) with a new value
#  - Column A widened
#  - Selection left on D3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: re-enter the phone number as text (leading apostrophe forces a
# text/quote-prefixed cell, matching the quotePrefix style Excel applies)
# before touching A2, so the new numeric-looking string gets its own
# shared-string slot ahead of the edited name text.
$ws.Range("B2").Value = "'792121074"

# A2: updated name text
$ws.Range("A2").Value = "Fayoumi11 Test Test2 Automation 51"

# Widen column A
$ws.Columns.Item(1).ColumnWidth = 32

# Leave the selection on D3, as in the saved file
$ws.Range("D3").Select() | Out-Null
